# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect a refreshed data scrape (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - row => new F value
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value  = 233
$wsExhibit.Range("F5").Value  = 6738
$wsExhibit.Range("F7").Value  = 434
$wsExhibit.Range("F8").Value  = 142
$wsExhibit.Range("F9").Value  = 6268
$wsExhibit.Range("F12").Value = 1264
$wsExhibit.Range("F21").Value = 4609
$wsExhibit.Range("F24").Value = 59
$wsExhibit.Range("F26").Value = 72

# Sheet "全部类型" (fourth sheet) - same updates, shifted by one row for the
# last two entries (F25/F27 instead of F24/F26)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value  = 233
$wsAll.Range("F5").Value  = 6738
$wsAll.Range("F7").Value  = 434
$wsAll.Range("F8").Value  = 142
$wsAll.Range("F9").Value  = 6268
$wsAll.Range("F12").Value = 1264
$wsAll.Range("F21").Value = 4609
$wsAll.Range("F25").Value = 59
$wsAll.Range("F27").Value = 72
